# culture_collection を MIxS から再度削除 (INSDC2017 での確認に基づく)
# Remove the "culture_collection" column (column AA) from the sheet.
# This shifts all columns after AA one position to the left, which affects
# both the worksheet data (via the shared-string table) and the per-column
# comments placed on header row 15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$deletedColumn = 27   # column AA

# 1. Snapshot every existing comment (row, column, text) before we touch
#    anything, because deleting/shifting cells does not automatically move
#    the attached comments in this environment.
$comments = $ws.Comments
$n = $comments.Count
$rows = New-Object 'object[]' $n
$cols = New-Object 'object[]' $n
$texts = New-Object 'object[]' $n
for ($i = 1; $i -le $n; $i++) {
    $c = $comments.Item($i)
    $rows[$i - 1] = $c.Parent.Row
    $cols[$i - 1] = $c.Parent.Column
    $texts[$i - 1] = $c.Text()
}

# 2. Remove all existing comments so we can re-create them at their
#    shifted locations afterwards.
for ($i = 1; $i -le $n; $i++) {
    $comments.Item(1).Delete()
}

# 3. Delete the culture_collection column itself, shifting everything to
#    its right one column to the left.
$ws.Columns(27).Delete()

# 4. Re-create the comments at their new column positions: comments that
#    were in the deleted column are dropped, comments to the right of it
#    move one column to the left, comments to the left are unaffected.
for ($i = 0; $i -lt $n; $i++) {
    $r = $rows[$i]
    $c = $cols[$i]
    $t = $texts[$i]

    if ($c -eq $deletedColumn) {
        continue
    }

    if ($c -gt $deletedColumn) {
        $c = $c - 1
    }

    $cell = $ws.Cells.Item($r, $c)
    $cell.ClearComments()
    $cell.AddComment($t)
}
